$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.944.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.733.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.76%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.732.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.59%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +2.20%  "

$ws.Range("E10").Value = "  +4.19%  "

$ws.Range("E11").Value = "  +2.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.357.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.732.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.876.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.80%  "

$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "497.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.73%  "

$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("E24").Value = "  +1.22%  "

$ws.Range("E25").Value = "  +1.29%  "

$ws.Range("E26").Value = "  -2.16%  "

$ws.Range("E27").Value = "  +0.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.875.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("E35").Value = "  +0.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.660.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.23%  "

$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "436.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("E45").Value = "  +1.86%  "

$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.62%  "

$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.742.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.45%  "
